# "Edit the Skill share" - adds a new "EditShareSkill" worksheet (a copy/variant
# of the SkillShareAdd form used to edit an existing skill share) and tweaks a
# couple of view-selection details on the existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: the engine stores column widths on a 1/6-character grid (stored =
# round(chars*6+5)/6). Given a desired *stored* width (as produced by Excel's
# real AutoFit/bestFit), back-solve the ColumnWidth to feed the COM setter so
# the round-tripped value lands as close as possible to the target.
# ---------------------------------------------------------------------------
function ToColumnWidth($storedWidth) {
    return ($storedWidth * 6 - 5) / 6
}

# ---------------------------------------------------------------------------
# 1. SkillShareAdd (sheet3): selection / scroll tweaks only.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SkillShareAdd")
$ws3.Activate()
$ws3.Range("F1").Select()
$ws3.Range("B1:Q2").Select()

$ws3.Columns.Item(18).ColumnWidth = ToColumnWidth 12.28515625
$ws3.Columns.Item(19).ColumnWidth = ToColumnWidth 15.7109375

# ---------------------------------------------------------------------------
# 2. Profile (sheet4): selection tweak only.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Profile")
$ws4.Range("E18").Select()

# ---------------------------------------------------------------------------
# 3. New sheet "EditShareSkill" appended after the last sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "EditShareSkill"

# Header row (row 1) - plain strings, no styling.
$ws5.Range("A1").Value = "OldTitle"
$ws5.Range("B1").Value = "OldDescription"
$ws5.Range("C1").Value = "Title"
$ws5.Range("D1").Value = "Description"
$ws5.Range("E1").Value = "Category"
$ws5.Range("F1").Value = "Subcategory"
$ws5.Range("G1").Value = "Tags"
$ws5.Range("H1").Value = "ServiceType"
$ws5.Range("I1").Value = "LocationType"
$ws5.Range("J1").Value = "StartDate"
$ws5.Range("K1").Value = "EndDate"
$ws5.Range("L1").Value = "DayOneSel"
$ws5.Range("M1").Value = "DayTwoSel"
$ws5.Range("N1").Value = "DayThreeSel"
$ws5.Range("O1").Value = "SkillTrade"
$ws5.Range("P1").Value = "Credit"
$ws5.Range("Q1").Value = "SkillExchange"
$ws5.Range("R1").Value = "Work Samples"
$ws5.Range("S1").Value = "Active"

# Data row (row 2) - left-aligned ("s=3"); set alignment before values/format
# so the engine reuses the pre-existing style slots instead of minting a new,
# unused one.
$ws5.Range("A2:S2").HorizontalAlignment = -4131

# Value insertion order matches the original authoring session so new shared
# strings land at the same indices as the reference workbook.
$ws5.Range("A2").Value = "Skill share"
$ws5.Range("B2").Value = "Basic"
$ws5.Range("Q2").Value = "test2"
$ws5.Range("C2").Value = "New Skill Share"
$ws5.Range("D2").Value = "I love automation testing"
$ws5.Range("E2").Value = "Graphics & Design"
$ws5.Range("F2").Value = "Flyers & Brochures"
$ws5.Range("G2").Value = "testSkills"
$ws5.Range("H2").Value = "One-off service"
$ws5.Range("I2").Value = "On-Site"

$ws5.Range("J2").NumberFormat = "d/mm/yyyy;@"
$ws5.Range("K2").NumberFormat = "d/mm/yyyy;@"
$ws5.Range("J2").Value = 43534
$ws5.Range("K2").Value = 43554

$ws5.Range("L2").Value = "Mon,01:00AM,03:00PM"
$ws5.Range("M2").Value = "Wed,01:00AM,03:00PM"
$ws5.Range("N2").Value = "Fri,01:00AM,03:00PM"
$ws5.Range("O2").Value = "Credit"
$ws5.Range("P2").Value = 10
$ws5.Range("R2").Value = "FileUpload.exe"
$ws5.Range("S2").Value = "Active"

# Row 3 stub: only E3/F3 exist, left-aligned, no value.
$ws5.Range("E3:F3").HorizontalAlignment = -4131

# Column widths (best effort - see ToColumnWidth note above).
$ws5.Columns.Item(1).ColumnWidth = ToColumnWidth 10
$ws5.Columns.Item(2).ColumnWidth = ToColumnWidth 14.42578125
$ws5.Columns.Item(3).ColumnWidth = ToColumnWidth 14.85546875
$ws5.Columns.Item(4).ColumnWidth = ToColumnWidth 23.5703125
$ws5.Columns.Item(5).ColumnWidth = ToColumnWidth 18.140625
$ws5.Columns.Item(6).ColumnWidth = ToColumnWidth 19.28515625
$ws5.Columns.Item(7).ColumnWidth = ToColumnWidth 9
$ws5.Columns.Item(8).ColumnWidth = ToColumnWidth 14.85546875
$ws5.Columns.Item(9).ColumnWidth = ToColumnWidth 12.7109375
$ws5.Columns.Item(10).ColumnWidth = ToColumnWidth 10.7109375
$ws5.Columns.Item(11).ColumnWidth = ToColumnWidth 10.7109375
$ws5.Columns.Item(12).ColumnWidth = ToColumnWidth 21.5703125
$ws5.Columns.Item(13).ColumnWidth = ToColumnWidth 21.7109375
$ws5.Columns.Item(14).ColumnWidth = ToColumnWidth 19.7109375
$ws5.Columns.Item(15).ColumnWidth = ToColumnWidth 9.7109375
$ws5.Columns.Item(16).ColumnWidth = ToColumnWidth 6.42578125
$ws5.Columns.Item(17).ColumnWidth = ToColumnWidth 13.140625
$ws5.Columns.Item(18).ColumnWidth = ToColumnWidth 14.5703125

# Final view state for the new sheet: F11 selected, tab active.
$ws5.Activate()
$ws5.Range("F11").Select()
